# Update "想去人数" (number of people interested) figures that changed
# between the two data refreshes, on both the "展览" sheet and the
# combined "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12330
$ws1.Range("F3").Value = 259
$ws1.Range("F4").Value = 250
$ws1.Range("F7").Value = 12268
$ws1.Range("F10").Value = 122
$ws1.Range("F13").Value = 6008
$ws1.Range("F15").Value = 3579

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12330
$ws4.Range("F3").Value = 259
$ws4.Range("F4").Value = 250
$ws4.Range("F8").Value = 12268
$ws4.Range("F11").Value = 122
$ws4.Range("F15").Value = 6008
$ws4.Range("F17").Value = 3579
